$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 10.10611590299825
$ws.Range("D2").Value = 5.095328314354823
$ws.Range("E2").Value = 13.42746585343532
$ws.Range("F2").Value = 25.6518663811093
$ws.Range("G2").Value = 3.637794530770863
$ws.Range("L2").Value = 9.736333479142679
$ws.Range("M2").Value = 59.49557385359812
$ws.Range("O2").Value = 22.91392056731097
$ws.Range("C3").Value = 10.31669980225114
$ws.Range("D3").Value = 5.143663013515512
$ws.Range("E3").Value = 13.21881991619042
$ws.Range("F3").Value = 25.97515904604353
$ws.Range("G3").Value = 3.642440724711932
$ws.Range("L3").Value = 9.56128167891673
$ws.Range("M3").Value = 56.10111810328626
$ws.Range("O3").Value = 23.26526068307782
$ws.Range("C4").Value = 10.45264581018219
$ws.Range("D4").Value = 5.175064505166724
$ws.Range("E4").Value = 13.09668431511862
$ws.Range("F4").Value = 26.19242375369575
$ws.Range("G4").Value = 3.645410858064896
$ws.Range("L4").Value = 9.45570469320057
$ws.Range("M4").Value = 53.90166570820146
$ws.Range("O4").Value = 23.49536023483869
$ws.Range("C5").Value = 10.50969399564524
$ws.Range("D5").Value = 5.188290963917646
$ws.Range("E5").Value = 13.04846354268325
$ws.Range("F5").Value = 26.28556321486365
$ws.Range("G5").Value = 3.646650933089758
$ws.Range("L5").Value = 9.413211889256633
$ws.Range("M5").Value = 52.97667278962672
$ws.Range("O5").Value = 23.59267711251196
$ws.Range("C6").Value = 10.51926580454994
$ws.Range("D6").Value = 5.190513091706482
$ws.Range("E6").Value = 13.04055151087475
$ws.Range("F6").Value = 26.3013036269783
$ws.Range("G6").Value = 3.646858647826607
$ws.Range("L6").Value = 9.406189549903111
$ws.Range("M6").Value = 52.8213550370209
$ws.Range("O6").Value = 23.6090490481246
$ws.Range("C7").Value = 10.45340853227658
$ws.Range("D7").Value = 5.175241144573855
$ws.Range("E7").Value = 13.09602765368839
$ws.Range("F7").Value = 26.19366137971527
$ws.Range("G7").Value = 3.645427461530641
$ws.Range("L7").Value = 9.45512940354805
$ws.Range("M7").Value = 53.8893066946526
$ws.Range("O7").Value = 23.49665839609424
$ws.Range("C8").Value = 10.17733147042538
$ws.Range("D8").Value = 5.111634226130922
$ws.Range("E8").Value = 13.35431478205575
$ws.Range("F8").Value = 25.75937294994932
$ws.Range("G8").Value = 3.639372312387642
$ws.Range("L8").Value = 9.675610944166548
$ws.Range("M8").Value = 58.34922674327721
$ws.Range("O8").Value = 23.0320385207732
$ws.Range("C9").Value = 9.689833071102576
$ws.Range("D9").Value = 5.000745505084448
$ws.Range("E9").Value = 13.90596424581097
$ws.Range("F9").Value = 25.06198192284597
$ws.Range("G9").Value = 3.62841926270772
$ws.Range("L9").Value = 10.12086103141774
$ws.Range("M9").Value = 66.17238254153467
$ws.Range("O9").Value = 22.23811643039232
$ws.Range("C10").Value = 9.366368961533562
$ws.Range("D10").Value = 4.927971007705028
$ws.Range("E10").Value = 14.33566642322749
$ws.Range("F10").Value = 24.65134054368261
$ws.Range("G10").Value = 3.620918963936241
$ws.Range("L10").Value = 10.4529409938505
$ws.Range("M10").Value = 71.35089414462129
$ws.Range("O10").Value = 21.73091758417123
$ws.Range("C11").Value = 9.22724726258541
$ws.Range("D11").Value = 4.896816649848311
$ws.Range("E11").Value = 14.53577602565821
$ws.Range("F11").Value = 24.48837128584228
$ws.Range("G11").Value = 3.617622336802807
$ws.Range("L11").Value = 10.60451053187614
$ws.Range("M11").Value = 73.58274324506252
$ws.Range("O11").Value = 21.5178230871501
$ws.Range("C12").Value = 9.175768608925388
$ws.Range("D12").Value = 4.885306025052649
$ws.Range("E12").Value = 14.61216342968768
$ws.Range("F12").Value = 24.43023748965814
$ws.Range("G12").Value = 3.616390297626789
$ws.Range("L12").Value = 10.66193344514235
$ws.Range("M12").Value = 74.41008167231267
$ws.Range("O12").Value = 21.43976746239898
$ws.Range("C13").Value = 9.186801055810433
$ws.Range("D13").Value = 4.887772179720236
$ws.Range("E13").Value = 14.59568577339803
$ws.Range("F13").Value = 24.44259599414707
$ws.Range("G13").Value = 3.616654917579418
$ws.Range("L13").Value = 10.64956589803744
$ws.Range("M13").Value = 74.23269153409206
$ws.Range("O13").Value = 21.4564590646833
$ws.Range("C14").Value = 9.222987648901386
$ws.Range("D14").Value = 4.895863876777773
$ws.Range("E14").Value = 14.54204854997918
$ws.Range("F14").Value = 24.48351596667401
$ws.Range("G14").Value = 3.617520650590103
$ws.Range("L14").Value = 10.60923446381947
$ws.Range("M14").Value = 73.65116579970345
$ws.Range("O14").Value = 21.51134790559784
$ws.Range("C15").Value = 9.245311336416956
$ws.Range("D15").Value = 4.900857828310388
$ws.Range("E15").Value = 14.50927202946939
$ws.Range("F15").Value = 24.50905132357781
$ws.Range("G15").Value = 3.618053054861376
$ws.Range("L15").Value = 10.58453247983229
$ws.Range("M15").Value = 73.29264503081721
$ws.Range("O15").Value = 21.54531571673626
$ws.Range("C16").Value = 9.375626466332051
$ws.Range("D16").Value = 4.930046803813864
$ws.Range("E16").Value = 14.32267721738075
$ws.Range("F16").Value = 24.6624834399929
$ws.Range("G16").Value = 3.621136704521297
$ws.Range("L16").Value = 10.44304174723973
$ws.Range("M16").Value = 71.20254402110731
$ws.Range("O16").Value = 21.74520727290193
$ws.Range("C17").Value = 9.457659107660977
$ws.Range("D17").Value = 4.948457313251619
$ws.Range("E17").Value = 14.20935408443424
$ws.Range("F17").Value = 24.76282003361363
$ws.Range("G17").Value = 3.623057773457798
$ws.Range("L17").Value = 10.35633757893966
$ws.Range("M17").Value = 69.88859565467529
$ws.Range("O17").Value = 21.8724215083762
$ws.Range("C18").Value = 9.505593978676249
$ws.Range("D18").Value = 4.959229789761626
$ws.Range("E18").Value = 14.14461326832235
$ws.Range("F18").Value = 24.82276493411979
$ws.Range("G18").Value = 3.624173589683288
$ws.Range("L18").Value = 10.30651770608612
$ws.Range("M18").Value = 69.12118485935399
$ws.Range("O18").Value = 21.94724485503059
$ws.Range("C19").Value = 9.521951621861232
$ws.Range("D19").Value = 4.962908454585035
$ws.Range("E19").Value = 14.12277039627248
$ws.Range("F19").Value = 24.84344037666553
$ws.Range("G19").Value = 3.6245532602062
$ws.Range("L19").Value = 10.289659593714
$ws.Range("M19").Value = 68.85935022239195
$ws.Range("O19").Value = 21.97285965049631
$ws.Range("C20").Value = 9.448848437218246
$ws.Range("D20").Value = 4.946478474953206
$ws.Range("E20").Value = 14.22137243700314
$ws.Range("F20").Value = 24.75190674132723
$ws.Range("G20").Value = 3.622852149412373
$ws.Range("L20").Value = 10.36556254404256
$ws.Range("M20").Value = 70.0296748173248
$ws.Range("O20").Value = 21.85870752452566
$ws.Range("C21").Value = 9.212325669622958
$ws.Range("D21").Value = 4.893479310991101
$ws.Range("E21").Value = 14.5577869813727
$ws.Range("F21").Value = 24.4713984095335
$ws.Range("G21").Value = 3.61726592294143
$ws.Range("L21").Value = 10.62108040601348
$ws.Range("M21").Value = 73.82245733935962
$ws.Range("O21").Value = 21.49515325294007
$ws.Range("C22").Value = 9.064788648291268
$ws.Range("D22").Value = 4.860516647481975
$ws.Range("E22").Value = 14.78118555233986
$ws.Range("F22").Value = 24.30900544350633
$ws.Range("G22").Value = 3.613710011979133
$ws.Range("L22").Value = 10.78821331032586
$ws.Range("M22").Value = 76.19744360099241
$ws.Range("O22").Value = 21.27298260761975
$ws.Range("C23").Value = 9.14286910664768
$ws.Range("D23").Value = 4.877953913853954
$ws.Range("E23").Value = 14.66164846404137
$ws.Range("F23").Value = 24.39371070631758
$ws.Range("G23").Value = 3.615599262508329
$ws.Range("L23").Value = 10.69901312984888
$ws.Range("M23").Value = 74.93936069436391
$ws.Range("O23").Value = 21.39011116889658
$ws.Range("C24").Value = 9.452829336630568
$ws.Range("D24").Value = 4.947372522546486
$ws.Range("E24").Value = 14.21593765793912
$ws.Range("F24").Value = 24.75683361264927
$ws.Range("G24").Value = 3.622945076615284
$ws.Range("L24").Value = 10.36139184813132
$ws.Range("M24").Value = 69.96593031850452
$ws.Range("O24").Value = 21.86490237223584
$ws.Range("C25").Value = 9.815770661416581
$ws.Range("D25").Value = 5.029236999251143
$ws.Range("E25").Value = 13.75221297396007
$ws.Range("F25").Value = 25.2333552949198
$ws.Range("G25").Value = 3.631285226717909
$ws.Range("L25").Value = 9.999335104676701
$ws.Range("M25").Value = 64.15568675200575
$ws.Range("O25").Value = 22.43988783682022
